$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 26667.846
$ws.Range("I2").Value = 34622.5
$ws.Range("J2").Value = 21696.188
$ws.Range("K2").Value = 34622.5
$ws.Range("L2").Value = 21696.188
$ws.Range("M2").Value = -34509.5
$ws.Range("N2").Value = -21922.188
$ws.Range("H9").Value = 383.9
$ws.Range("I9").Value = 252.875
$ws.Range("J9").Value = 908
$ws.Range("K9").Value = 252.875
$ws.Range("L9").Value = 908
$ws.Range("M9").Value = -83.875
$ws.Range("N9").Value = -1246
$ws.Range("H28").Value = 851.9048
$ws.Range("J28").Value = 293.33334
$ws.Range("L28").Value = 293.33334
$ws.Range("N28").Value = -1263.33334
$ws.Range("H103").Value = 909.8
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 909.8
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2729.4
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -3901.4
$ws.Range("H138").Value = 3691.81
$ws.Range("J138").Value = 4016.8572
$ws.Range("L138").Value = 12050.5716
$ws.Range("N138").Value = -22330.5716
$ws.Range("H139").Value = 289990.94
$ws.Range("J139").Value = 289990.94
$ws.Range("L139").Value = 289990.94
$ws.Range("N139").Value = -300270.94
$ws.Range("H141").Value = 3214.5715
$ws.Range("I141").Value = 3178.2222
$ws.Range("K141").Value = 9534.6666
$ws.Range("M141").Value = -4354.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2980.0637
$ws.Range("I61").Value = 1801.8158
$ws.Range("J61").Value = 7954.8887
$ws.Range("K61").Value = 1801.8158
$ws.Range("L61").Value = 7954.8887
$ws.Range("M61").Value = -1589.8158
$ws.Range("N61").Value = -8378.8887
$ws.Range("H110").Value = 12809.117
$ws.Range("I110").Value = 13462.167
$ws.Range("K110").Value = 13462.167
$ws.Range("M110").Value = -11417.167
$ws.Range("H136").Value = 2980.0637
$ws.Range("I136").Value = 1801.8158
$ws.Range("J136").Value = 7954.8887
$ws.Range("K136").Value = 5405.4474
$ws.Range("L136").Value = 23864.6661
$ws.Range("M136").Value = -2855.4474
$ws.Range("N136").Value = -28964.6661
$ws.Range("H139").Value = 84166.336
$ws.Range("J139").Value = 84166.336
$ws.Range("L139").Value = 84166.336
$ws.Range("N139").Value = -94446.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1643.0588
$ws.Range("I107").Value = 1247.8462
$ws.Range("J107").Value = 2927.5
$ws.Range("K107").Value = 1247.8462
$ws.Range("L107").Value = 2927.5
$ws.Range("M107").Value = 672.1538
$ws.Range("N107").Value = -6767.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30155.566
$ws.Range("I31").Value = 3542.25
$ws.Range("J31").Value = 47897.777
$ws.Range("K31").Value = 3542.25
$ws.Range("L31").Value = 47897.777
$ws.Range("M31").Value = -3247.25
$ws.Range("N31").Value = -48487.777
$ws.Range("H34").Value = 30155.566
$ws.Range("I34").Value = 3542.25
$ws.Range("J34").Value = 47897.777
$ws.Range("K34").Value = 3542.25
$ws.Range("L34").Value = 47897.777
$ws.Range("M34").Value = -3340.25
$ws.Range("N34").Value = -48301.777
$ws.Range("H58").Value = 3860.3635
$ws.Range("I58").Value = 3327.5
$ws.Range("K58").Value = 3327.5
$ws.Range("M58").Value = -3124.5
$ws.Range("H122").Value = 3229.1333
$ws.Range("J122").Value = 3554
$ws.Range("L122").Value = 10662
$ws.Range("N122").Value = -15562
$ws.Range("H132").Value = 2911.2307
$ws.Range("I132").Value = 3059.3157
$ws.Range("K132").Value = 9177.947100000001
$ws.Range("M132").Value = -6647.947100000001
$ws.Range("H134").Value = 4628.972
$ws.Range("I134").Value = 5068.7744
$ws.Range("J134").Value = 1902.2
$ws.Range("K134").Value = 15206.3232
$ws.Range("L134").Value = 5706.6
$ws.Range("M134").Value = -12671.3232
$ws.Range("N134").Value = -10776.6
$ws.Range("H136").Value = 3860.3635
$ws.Range("I136").Value = 3327.5
$ws.Range("K136").Value = 9982.5
$ws.Range("M136").Value = -7432.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8334270
$ws.Range("I68").Value = 16667465
$ws.Range("J68").Value = 1074.8
$ws.Range("K68").Value = 50002395
$ws.Range("L68").Value = 3224.4
$ws.Range("M68").Value = -50001584
$ws.Range("N68").Value = -4846.4
$ws.Range("H71").Value = 8334270
$ws.Range("I71").Value = 16667465
$ws.Range("J71").Value = 1074.8
$ws.Range("K71").Value = 150007185
$ws.Range("L71").Value = 9673.199999999999
$ws.Range("M71").Value = -150003129
$ws.Range("N71").Value = -17785.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 19792
$ws.Range("J48").Value = 19792
$ws.Range("L48").Value = 19792
$ws.Range("N48").Value = -20762
$ws.Range("H113").Value = 22504.125
$ws.Range("I113").Value = 28601.6
$ws.Range("K113").Value = 28601.6
$ws.Range("M113").Value = -26431.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1236.8334
$ws.Range("I46").Value = 1236.8334
$ws.Range("K46").Value = 1236.8334
$ws.Range("M46").Value = -1048.8334
$ws.Range("H61").Value = 1667.8422
$ws.Range("I61").Value = 1667.8422
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1667.8422
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1465.8422
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1667.8422
$ws.Range("I113").Value = 1667.8422
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1667.8422
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 502.1578
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 8695.666999999999
$ws.Range("I122").Value = 8579.076999999999
$ws.Range("K122").Value = 25737.231
$ws.Range("M122").Value = -23287.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20843648
$ws.Range("I62").Value = 6743
$ws.Range("K62").Value = 6743
$ws.Range("M62").Value = -6119
$ws.Range("H65").Value = 20843648
$ws.Range("I65").Value = 6743
$ws.Range("K65").Value = 33715
$ws.Range("M65").Value = -30595
$ws.Range("H96").Value = 69682
$ws.Range("J96").Value = 4988.25
$ws.Range("L96").Value = 4988.25
$ws.Range("N96").Value = -7734.25
$ws.Range("H112").Value = 24908.334
$ws.Range("J112").Value = 25400
$ws.Range("L112").Value = 25400
$ws.Range("N112").Value = -28354
$ws.Range("H113").Value = 1741.1724
$ws.Range("J113").Value = 1835.2727
$ws.Range("L113").Value = 5505.8181
$ws.Range("N113").Value = -9845.8181
$ws.Range("H136").Value = 3531.8738
$ws.Range("I136").Value = 3206.5322
$ws.Range("J136").Value = 4143.121
$ws.Range("K136").Value = 9619.596600000001
$ws.Range("L136").Value = 12429.363
$ws.Range("M136").Value = -7069.596600000001
$ws.Range("N136").Value = -17529.363
